# Auto-generated edit script applying the diff's cell-level numeric updates.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2067.1765
$ws.Range("I112").Value = 1249.5
$ws.Range("J112").Value = 2318.7693
$ws.Range("K112").Value = 3748.5
$ws.Range("L112").Value = 6956.3079
$ws.Range("M112").Value = -2640.5
$ws.Range("N112").Value = -9172.3079
$ws.Range("H135").Value = 2380.2593
$ws.Range("I135").Value = 2292.125
$ws.Range("J135").Value = 2508.4546
$ws.Range("K135").Value = 20629.125
$ws.Range("L135").Value = 22576.0914
$ws.Range("M135").Value = -18094.125
$ws.Range("N135").Value = -27646.0914

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H40").Value = 7342.6665
$ws.Range("J40").Value = 6000
$ws.Range("L40").Value = 6000
$ws.Range("N40").Value = -6352
$ws.Range("H122").Value = 6602
$ws.Range("I122").Value = 6652.25
$ws.Range("J122").Value = 6501.5
$ws.Range("K122").Value = 19956.75
$ws.Range("L122").Value = 19504.5
$ws.Range("M122").Value = -17506.75
$ws.Range("N122").Value = -24404.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 15848.625
$ws.Range("J95").Value = 15848.625
$ws.Range("L95").Value = 15848.625
$ws.Range("N95").Value = -21340.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 18499.857
$ws.Range("I55").Value = 6125
$ws.Range("J55").Value = 34999.668
$ws.Range("K55").Value = 6125
$ws.Range("L55").Value = 34999.668
$ws.Range("M55").Value = -5810
$ws.Range("N55").Value = -35629.668
$ws.Range("H58").Value = 4680.294
$ws.Range("I58").Value = 2955.5557
$ws.Range("K58").Value = 2955.5557
$ws.Range("M58").Value = -2752.5557
$ws.Range("H105").Value = 3602.5
$ws.Range("J105").Value = 6000
$ws.Range("L105").Value = 6000
$ws.Range("N105").Value = -9494
$ws.Range("H132").Value = 16413.441
$ws.Range("I132").Value = 5151.8335
$ws.Range("K132").Value = 15455.5005
$ws.Range("M132").Value = -12925.5005
$ws.Range("H134").Value = 6388.12
$ws.Range("I134").Value = 6397
$ws.Range("K134").Value = 19191
$ws.Range("M134").Value = -16656
$ws.Range("H136").Value = 4680.294
$ws.Range("I136").Value = 2955.5557
$ws.Range("K136").Value = 8866.667099999999
$ws.Range("M136").Value = -6316.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2371.476
$ws.Range("I34").Value = 919.8889
$ws.Range("J34").Value = 3460.1667
$ws.Range("K34").Value = 2759.6667
$ws.Range("L34").Value = 10380.5001
$ws.Range("M34").Value = -2675.6667
$ws.Range("N34").Value = -10548.5001
$ws.Range("H36").Value = 496
$ws.Range("I36").Value = 496
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1488
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1319
$ws.Range("N36").Value = ""
$ws.Range("H38").Value = 66666684
$ws.Range("J38").Value = 17.363636
$ws.Range("L38").Value = 52.090908
$ws.Range("N38").Value = -746.090908
$ws.Range("H39").Value = 3919.625
$ws.Range("I39").Value = 875
$ws.Range("J39").Value = 6964.25
$ws.Range("K39").Value = 2625
$ws.Range("L39").Value = 20892.75
$ws.Range("M39").Value = -2331
$ws.Range("N39").Value = -21480.75
$ws.Range("H55").Value = 1385.2106
$ws.Range("I55").Value = 741.6667
$ws.Range("J55").Value = 1682.2307
$ws.Range("K55").Value = 2225.0001
$ws.Range("L55").Value = 5046.6921
$ws.Range("M55").Value = -2048.0001
$ws.Range("N55").Value = -5400.6921
$ws.Range("H56").Value = 12087
$ws.Range("I56").Value = 12087
$ws.Range("K56").Value = 12087
$ws.Range("M56").Value = -11557
$ws.Range("H57").Value = 2889.2222
$ws.Range("I57").Value = 500
$ws.Range("J57").Value = 3571.8572
$ws.Range("K57").Value = 1500
$ws.Range("L57").Value = 10715.5716
$ws.Range("M57").Value = -941
$ws.Range("N57").Value = -11833.5716
$ws.Range("H68").Value = 6933.3335
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 9900
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 29700
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -31322
$ws.Range("H71").Value = 6933.3335
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 9900
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 89100
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -97212
$ws.Range("H80").Value = 3000.25
$ws.Range("H83").Value = 3000.25
$ws.Range("H122").Value = 55560480
$ws.Range("I122").Value = 111119896
$ws.Range("J122").Value = 1067
$ws.Range("K122").Value = 1000079064
$ws.Range("L122").Value = 9603
$ws.Range("M122").Value = -1000076614
$ws.Range("N122").Value = -14503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10613.714
$ws.Range("I70").Value = 8574.75
$ws.Range("K70").Value = 8574.75
$ws.Range("M70").Value = -8304.75
$ws.Range("H73").Value = 10613.714
$ws.Range("I73").Value = 8574.75
$ws.Range("K73").Value = 8574.75
$ws.Range("M73").Value = -7638.75
$ws.Range("H136").Value = 30162.96
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 30162.96
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 90488.88
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -95588.88

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 46214.2
$ws.Range("I74").Value = 39707.332
$ws.Range("K74").Value = 39707.332
$ws.Range("M74").Value = -38709.332
$ws.Range("H77").Value = 46214.2
$ws.Range("I77").Value = 39707.332
$ws.Range("K77").Value = 119121.996
$ws.Range("M77").Value = -114129.996
$ws.Range("H132").Value = 4973.161
$ws.Range("I132").Value = 4190.548
$ws.Range("K132").Value = 12571.644
$ws.Range("M132").Value = -10041.644

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1161.75
$ws.Range("I107").Value = 1065.3871
$ws.Range("K107").Value = 3196.1613
$ws.Range("M107").Value = -1276.1613
$ws.Range("H123").Value = 49999.5
$ws.Range("J123").Value = 49999.5
$ws.Range("L123").Value = 49999.5
$ws.Range("N123").Value = -59799.5
$ws.Range("H126").Value = 8678.306
$ws.Range("I126").Value = 7237
$ws.Range("K126").Value = 21711
$ws.Range("M126").Value = -19241
